$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1000
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H64").Value = 6500
$ws.Range("I64").Value = 5000
$ws.Range("K64").Value = 5000
$ws.Range("M64").Value = -4752
$ws.Range("H67").Value = 6500
$ws.Range("I67").Value = 5000
$ws.Range("K67").Value = 5000
$ws.Range("M67").Value = -4142
$ws.Range("H86").Value = 4052160.8
$ws.Range("I86").Value = 3468.8333
$ws.Range("J86").Value = 7522468
$ws.Range("K86").Value = 3468.8333
$ws.Range("L86").Value = 7522468
$ws.Range("M86").Value = -2345.8333
$ws.Range("N86").Value = -7524714
$ws.Range("H89").Value = 4052160.8
$ws.Range("I89").Value = 3468.8333
$ws.Range("J89").Value = 7522468
$ws.Range("K89").Value = 17344.1665
$ws.Range("L89").Value = 37612340
$ws.Range("M89").Value = -11728.1665
$ws.Range("N89").Value = -37623572
$ws.Range("H96").Value = 585.375
$ws.Range("I96").Value = 604.8570999999999
$ws.Range("K96").Value = 1814.5713
$ws.Range("M96").Value = -441.5712999999998
$ws.Range("H98").Value = 2371.6553
$ws.Range("I98").Value = 1217.2273
$ws.Range("K98").Value = 1217.2273
$ws.Range("M98").Value = 280.7727
$ws.Range("H122").Value = 2371.6553
$ws.Range("I122").Value = 1217.2273
$ws.Range("K122").Value = 3651.6819
$ws.Range("M122").Value = -1201.6819
$ws.Range("H126").Value = 73405.53
$ws.Range("J126").Value = 73405.53
$ws.Range("L126").Value = 73405.53
$ws.Range("N126").Value = -83285.53
$ws.Range("H127").Value = 11064.23
$ws.Range("I127").Value = 13303.8
$ws.Range("K127").Value = 39911.39999999999
$ws.Range("M127").Value = -34951.39999999999
$ws.Range("H128").Value = 88177.45
$ws.Range("J128").Value = 88177.45
$ws.Range("L128").Value = 88177.45
$ws.Range("N128").Value = -98137.45
$ws.Range("H129").Value = 13405.5
$ws.Range("I129").Value = 24160.4
$ws.Range("J129").Value = 9820.532999999999
$ws.Range("K129").Value = 72481.20000000001
$ws.Range("L129").Value = 29461.599
$ws.Range("M129").Value = -67481.20000000001
$ws.Range("N129").Value = -39461.599
$ws.Range("H131").Value = 3596.4285
$ws.Range("I131").Value = 3236.2
$ws.Range("K131").Value = 9708.599999999999
$ws.Range("M131").Value = -4668.599999999999
$ws.Range("H133").Value = 59548.54
$ws.Range("J133").Value = 59548.54
$ws.Range("L133").Value = 59548.54
$ws.Range("N133").Value = -69668.54000000001
$ws.Range("H137").Value = 4227.0464
$ws.Range("I137").Value = 2569.1614
$ws.Range("J137").Value = 8509.916999999999
$ws.Range("K137").Value = 7707.4842
$ws.Range("L137").Value = 25529.751
$ws.Range("M137").Value = -5157.4842
$ws.Range("N137").Value = -30629.751
$ws.Range("H138").Value = 6030.1816
$ws.Range("J138").Value = 6278.7188
$ws.Range("L138").Value = 18836.1564
$ws.Range("N138").Value = -29116.1564
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1943.4459
$ws.Range("I32").Value = 1409.058
$ws.Range("J32").Value = 9318
$ws.Range("K32").Value = 1409.058
$ws.Range("L32").Value = 9318
$ws.Range("M32").Value = -1122.058
$ws.Range("N32").Value = -9892
$ws.Range("H61").Value = 3494.25
$ws.Range("I61").Value = 1991.7858
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 1991.7858
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = -1779.7858
$ws.Range("N61").Value = -7424
$ws.Range("H74").Value = 1306.88
$ws.Range("I74").Value = 1070.1277
$ws.Range("J74").Value = 5016
$ws.Range("K74").Value = 1070.1277
$ws.Range("L74").Value = 5016
$ws.Range("M74").Value = -196.1277
$ws.Range("N74").Value = -6764
$ws.Range("H77").Value = 1306.88
$ws.Range("I77").Value = 1070.1277
$ws.Range("J77").Value = 5016
$ws.Range("K77").Value = 5350.6385
$ws.Range("L77").Value = 25080
$ws.Range("M77").Value = -982.6385
$ws.Range("N77").Value = -33816
$ws.Range("H97").Value = 1299.9375
$ws.Range("I97").Value = 1729.409
$ws.Range("J97").Value = 355.1
$ws.Range("K97").Value = 1729.409
$ws.Range("L97").Value = 355.1
$ws.Range("M97").Value = -1233.409
$ws.Range("N97").Value = -1347.1
$ws.Range("H102").Value = 1327.2858
$ws.Range("I102").Value = 1356.3235
$ws.Range("K102").Value = 1356.3235
$ws.Range("M102").Value = 265.6765
$ws.Range("H136").Value = 3494.25
$ws.Range("I136").Value = 1991.7858
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 5975.357400000001
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -3425.357400000001
$ws.Range("N136").Value = -26100

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2478.0881
$ws.Range("I31").Value = 1329.1364
$ws.Range("J31").Value = 4584.5
$ws.Range("K31").Value = 1329.1364
$ws.Range("L31").Value = 4584.5
$ws.Range("M31").Value = -1034.1364
$ws.Range("N31").Value = -5174.5
$ws.Range("H34").Value = 2478.0881
$ws.Range("I34").Value = 1329.1364
$ws.Range("J34").Value = 4584.5
$ws.Range("K34").Value = 1329.1364
$ws.Range("L34").Value = 4584.5
$ws.Range("M34").Value = -1127.1364
$ws.Range("N34").Value = -4988.5
$ws.Range("H86").Value = 63429.43
$ws.Range("I86").Value = 136335.33
$ws.Range("J86").Value = 8750
$ws.Range("K86").Value = 136335.33
$ws.Range("L86").Value = 8750
$ws.Range("M86").Value = -135212.33
$ws.Range("N86").Value = -10996
$ws.Range("H89").Value = 63429.43
$ws.Range("I89").Value = 136335.33
$ws.Range("J89").Value = 8750
$ws.Range("K89").Value = 681676.6499999999
$ws.Range("L89").Value = 43750
$ws.Range("M89").Value = -676060.6499999999
$ws.Range("N89").Value = -54982
$ws.Range("H105").Value = 1523.4286
$ws.Range("I105").Value = 1523.4286
$ws.Range("K105").Value = 1523.4286
$ws.Range("M105").Value = 223.5714

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 2974
$ws.Range("I43").Value = 2974
$ws.Range("K43").Value = 8922
$ws.Range("M43").Value = -8808
$ws.Range("H96").Value = 200406000
$ws.Range("I96").Value = 500500000
$ws.Range("J96").Value = 343333.34
$ws.Range("K96").Value = 1501500000
$ws.Range("L96").Value = 1030000.02
$ws.Range("M96").Value = -1501497941
$ws.Range("N96").Value = -1034118.02
$ws.Range("H107").Value = 56304.42
$ws.Range("I107").Value = 2369.6667
$ws.Range("J107").Value = 66417.19
$ws.Range("K107").Value = 7109.000100000001
$ws.Range("L107").Value = 199251.57
$ws.Range("M107").Value = -5189.000100000001
$ws.Range("N107").Value = -203091.57
$ws.Range("H113").Value = 2849383.2
$ws.Range("J113").Value = 455
$ws.Range("L113").Value = 1365
$ws.Range("N113").Value = -5705
$ws.Range("H141").Value = 31633
$ws.Range("I141").Value = 24949.5
$ws.Range("K141").Value = 74848.5
$ws.Range("M141").Value = -69668.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 8415.691999999999
$ws.Range("J97").Value = 508.33334
$ws.Range("L97").Value = 508.33334
$ws.Range("N97").Value = -1500.33334
$ws.Range("H102").Value = 1926.6428
$ws.Range("I102").Value = 2147.4
$ws.Range("J102").Value = 1374.75
$ws.Range("K102").Value = 2147.4
$ws.Range("L102").Value = 1374.75
$ws.Range("M102").Value = -525.4000000000001
$ws.Range("N102").Value = -4618.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("H46").Value = 3533.0688
$ws.Range("I46").Value = 3023.6843
$ws.Range("K46").Value = 3023.6843
$ws.Range("M46").Value = -2835.6843
$ws.Range("H132").Value = 5296.2
$ws.Range("I132").Value = 4206.2856
$ws.Range("J132").Value = 6249.875
$ws.Range("K132").Value = 12618.8568
$ws.Range("L132").Value = 18749.625
$ws.Range("M132").Value = -10088.8568
$ws.Range("N132").Value = -23809.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5213.857
$ws.Range("J126").Value = 6000
$ws.Range("L126").Value = 18000
$ws.Range("N126").Value = -22940
